$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Efna1"
$ws.Range("C2").Value = "Epha2"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 3.0
$ws.Range("F2").Value = 1.0
$ws.Range("G2").Value = 23.82411166666667
$ws.Range("H2").Value = 71.472335
$ws.Range("I2").Value = 0.8811513593020099
$ws.Range("J2").Value = 0.8811513593020099
$ws.Range("K2").Value = 3.0
$ws.Range("L2").Value = 1.0
$ws.Range("M2").Value = 24.244885
$ws.Range("N2").Value = 72.734655
$ws.Range("O2").Value = 0.6895205882382217
$ws.Range("P2").Value = 0.6895205882382218
$ws.Range("Q2").Value = 577.6128475854916
$ws.Range("R2").Value = 5198.515628269425
$ws.Range("S2").Value = 0.6075720035928306
$ws.Range("T2").Value = 0.6075720035928306

# Row 3
$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Efna1"
$ws.Range("C3").Value = "Epha2"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 3.0
$ws.Range("F3").Value = 1.0
$ws.Range("G3").Value = 23.82411166666667
$ws.Range("H3").Value = 71.472335
$ws.Range("I3").Value = 0.8811513593020099
$ws.Range("J3").Value = 0.8811513593020099
$ws.Range("K3").Value = 3.0
$ws.Range("L3").Value = 1.0
$ws.Range("M3").Value = 0.7651789999999999
$ws.Range("N3").Value = 2.295537
$ws.Range("O3").Value = 0.02176156637523891
$ws.Range("P3").Value = 0.02176156637523891
$ws.Range("Q3").Value = 18.22970994098833
$ws.Range("R3").Value = 164.067389468895
$ws.Range("S3").Value = 0.01917523379208268
$ws.Range("T3").Value = 0.01917523379208268

# Row 4
$ws.Range("A4").Value = "ECs"
$ws.Range("B4").Value = "Efna1"
$ws.Range("C4").Value = "Epha2"
$ws.Range("D4").Value = "M2"
$ws.Range("E4").Value = 3.0
$ws.Range("F4").Value = 1.0
$ws.Range("G4").Value = 23.82411166666667
$ws.Range("H4").Value = 71.472335
$ws.Range("I4").Value = 0.8811513593020099
$ws.Range("J4").Value = 0.8811513593020099
$ws.Range("K4").Value = 3.0
$ws.Range("L4").Value = 1.0
$ws.Range("M4").Value = 0.6351283333333334
$ws.Range("N4").Value = 1.905385
$ws.Range("O4").Value = 0.01806294655581008
$ws.Range("P4").Value = 0.01806294655581008
$ws.Range("Q4").Value = 15.13136833599722
$ws.Range("R4").Value = 136.182315023975
$ws.Range("S4").Value = 0.01591618991065161
$ws.Range("T4").Value = 0.01591618991065161

# Row 5
$ws.Range("A5").Value = "ECs"
$ws.Range("B5").Value = "Efna1"
$ws.Range("C5").Value = "Epha2"
$ws.Range("D5").Value = "sCs"
$ws.Range("E5").Value = 3.0
$ws.Range("F5").Value = 1.0
$ws.Range("G5").Value = 23.82411166666667
$ws.Range("H5").Value = 71.472335
$ws.Range("I5").Value = 0.8811513593020099
$ws.Range("J5").Value = 0.8811513593020099
$ws.Range("K5").Value = 3.0
$ws.Range("L5").Value = 1.0
$ws.Range("M5").Value = 9.516752666666667
$ws.Range("N5").Value = 28.550258
$ws.Range("O5").Value = 0.2706548988307292
$ws.Range("P5").Value = 0.2706548988307293
$ws.Range("Q5").Value = 226.7281782347144
$ws.Range("R5").Value = 2040.55360411243
$ws.Range("S5").Value = 0.2384879320064451
$ws.Range("T5").Value = 0.2384879320064451

# Row 6
$ws.Range("A6").Value = "FAPs"
$ws.Range("B6").Value = "Efna1"
$ws.Range("C6").Value = "Epha2"
$ws.Range("D6").Value = "ECs"
$ws.Range("E6").Value = 3.0
$ws.Range("F6").Value = 1.0
$ws.Range("G6").Value = 2.471212
$ws.Range("H6").Value = 7.413636
$ws.Range("I6").Value = 0.09139949658522162
$ws.Range("J6").Value = 0.09139949658522162
$ws.Range("K6").Value = 3.0
$ws.Range("L6").Value = 1.0
$ws.Range("M6").Value = 24.244885
$ws.Range("N6").Value = 72.734655
$ws.Range("O6").Value = 0.6895205882382217
$ws.Range("P6").Value = 0.6895205882382218
$ws.Range("Q6").Value = 59.91425075062
$ws.Range("R6").Value = 539.2282567555801
$ws.Range("S6").Value = 0.06302183465011935
$ws.Range("T6").Value = 0.06302183465011936

# Row 7
$ws.Range("A7").Value = "FAPs"
$ws.Range("B7").Value = "Efna1"
$ws.Range("C7").Value = "Epha2"
$ws.Range("D7").Value = "FAPs"
$ws.Range("E7").Value = 3.0
$ws.Range("F7").Value = 1.0
$ws.Range("G7").Value = 2.471212
$ws.Range("H7").Value = 7.413636
$ws.Range("I7").Value = 0.09139949658522162
$ws.Range("J7").Value = 0.09139949658522162
$ws.Range("K7").Value = 3.0
$ws.Range("L7").Value = 1.0
$ws.Range("M7").Value = 0.7651789999999999
$ws.Range("N7").Value = 2.295537
$ws.Range("O7").Value = 0.02176156637523891
$ws.Range("P7").Value = 0.02176156637523891
$ws.Range("Q7").Value = 1.890919526948
$ws.Range("R7").Value = 17.018275742532
$ws.Range("S7").Value = 0.001988996211602722
$ws.Range("T7").Value = 0.001988996211602723

# Row 8
$ws.Range("A8").Value = "FAPs"
$ws.Range("B8").Value = "Efna1"
$ws.Range("C8").Value = "Epha2"
$ws.Range("D8").Value = "M2"
$ws.Range("E8").Value = 3.0
$ws.Range("F8").Value = 1.0
$ws.Range("G8").Value = 2.471212
$ws.Range("H8").Value = 7.413636
$ws.Range("I8").Value = 0.09139949658522162
$ws.Range("J8").Value = 0.09139949658522162
$ws.Range("K8").Value = 3.0
$ws.Range("L8").Value = 1.0
$ws.Range("M8").Value = 0.6351283333333334
$ws.Range("N8").Value = 1.905385
$ws.Range("O8").Value = 0.01806294655581008
$ws.Range("P8").Value = 0.01806294655581008
$ws.Range("Q8").Value = 1.569536758873334
$ws.Range("R8").Value = 14.12583082986
$ws.Range("S8").Value = 0.001650944222046804
$ws.Range("T8").Value = 0.001650944222046804

# Row 9
$ws.Range("A9").Value = "FAPs"
$ws.Range("B9").Value = "Efna1"
$ws.Range("C9").Value = "Epha2"
$ws.Range("D9").Value = "sCs"
$ws.Range("E9").Value = 3.0
$ws.Range("F9").Value = 1.0
$ws.Range("G9").Value = 2.471212
$ws.Range("H9").Value = 7.413636
$ws.Range("I9").Value = 0.09139949658522162
$ws.Range("J9").Value = 0.09139949658522162
$ws.Range("K9").Value = 3.0
$ws.Range("L9").Value = 1.0
$ws.Range("M9").Value = 9.516752666666667
$ws.Range("N9").Value = 28.550258
$ws.Range("O9").Value = 0.2706548988307292
$ws.Range("P9").Value = 0.2706548988307293
$ws.Range("Q9").Value = 23.51791339089867
$ws.Range("R9").Value = 211.661220518088
$ws.Range("S9").Value = 0.02473772150145274
$ws.Range("T9").Value = 0.02473772150145274

# Row 10
$ws.Range("A10").Value = "M2"
$ws.Range("B10").Value = "Efna1"
$ws.Range("C10").Value = "Epha2"
$ws.Range("D10").Value = "ECs"
$ws.Range("E10").Value = 1.0
$ws.Range("F10").Value = 0.3333333333333333
$ws.Range("G10").Value = 0.02836866666666667
$ws.Range("H10").Value = 0.085106
$ws.Range("I10").Value = 0.001049234890461559
$ws.Range("J10").Value = 0.001049234890461559
$ws.Range("K10").Value = 3.0
$ws.Range("L10").Value = 1.0
$ws.Range("M10").Value = 24.244885
$ws.Range("N10").Value = 72.734655
$ws.Range("O10").Value = 0.6895205882382217
$ws.Range("P10").Value = 0.6895205882382218
$ws.Range("Q10").Value = 0.6877950609366666
$ws.Range("R10").Value = 6.190155548430001
$ws.Range("S10").Value = 0.0007234690588711203
$ws.Range("T10").Value = 0.0007234690588711205

# Row 11
$ws.Range("A11").Value = "M2"
$ws.Range("B11").Value = "Efna1"
$ws.Range("C11").Value = "Epha2"
$ws.Range("D11").Value = "FAPs"
$ws.Range("E11").Value = 1.0
$ws.Range("F11").Value = 0.3333333333333333
$ws.Range("G11").Value = 0.02836866666666667
$ws.Range("H11").Value = 0.085106
$ws.Range("I11").Value = 0.001049234890461559
$ws.Range("J11").Value = 0.001049234890461559
$ws.Range("K11").Value = 3.0
$ws.Range("L11").Value = 1.0
$ws.Range("M11").Value = 0.7651789999999999
$ws.Range("N11").Value = 2.295537
$ws.Range("O11").Value = 0.02176156637523891
$ws.Range("P11").Value = 0.02176156637523891
$ws.Range("Q11").Value = 0.02170710799133333
$ws.Range("R11").Value = 0.195363971922
$ws.Range("S11").Value = 0.00002283299471199574
$ws.Range("T11").Value = 0.00002283299471199575

# Row 12
$ws.Range("A12").Value = "M2"
$ws.Range("B12").Value = "Efna1"
$ws.Range("C12").Value = "Epha2"
$ws.Range("D12").Value = "M2"
$ws.Range("E12").Value = 1.0
$ws.Range("F12").Value = 0.3333333333333333
$ws.Range("G12").Value = 0.02836866666666667
$ws.Range("H12").Value = 0.085106
$ws.Range("I12").Value = 0.001049234890461559
$ws.Range("J12").Value = 0.001049234890461559
$ws.Range("K12").Value = 3.0
$ws.Range("L12").Value = 1.0
$ws.Range("M12").Value = 0.6351283333333334
$ws.Range("N12").Value = 1.905385
$ws.Range("O12").Value = 0.01806294655581008
$ws.Range("P12").Value = 0.01806294655581008
$ws.Range("Q12").Value = 0.01801774397888889
$ws.Range("R12").Value = 0.16215969581
$ws.Range("S12").Value = 0.00001895227375089838
$ws.Range("T12").Value = 0.00001895227375089838

# Row 13
$ws.Range("A13").Value = "M2"
$ws.Range("B13").Value = "Efna1"
$ws.Range("C13").Value = "Epha2"
$ws.Range("D13").Value = "sCs"
$ws.Range("E13").Value = 1.0
$ws.Range("F13").Value = 0.3333333333333333
$ws.Range("G13").Value = 0.02836866666666667
$ws.Range("H13").Value = 0.085106
$ws.Range("I13").Value = 0.001049234890461559
$ws.Range("J13").Value = 0.001049234890461559
$ws.Range("K13").Value = 3.0
$ws.Range("L13").Value = 1.0
$ws.Range("M13").Value = 9.516752666666667
$ws.Range("N13").Value = 28.550258
$ws.Range("O13").Value = 0.2706548988307292
$ws.Range("P13").Value = 0.2706548988307293
$ws.Range("Q13").Value = 0.2699775841497778
$ws.Range("R13").Value = 2.429798257348
$ws.Range("S13").Value = 0.0002839805631275446
$ws.Range("T13").Value = 0.0002839805631275446

# Row 14
$ws.Range("A14").Value = "sCs"
$ws.Range("B14").Value = "Efna1"
$ws.Range("C14").Value = "Epha2"
$ws.Range("D14").Value = "ECs"
$ws.Range("E14").Value = 3.0
$ws.Range("F14").Value = 1.0
$ws.Range("G14").Value = 0.713787
$ws.Range("H14").Value = 2.141361
$ws.Range("I14").Value = 0.02639990922230694
$ws.Range("J14").Value = 0.02639990922230694
$ws.Range("K14").Value = 3.0
$ws.Range("L14").Value = 1.0
$ws.Range("M14").Value = 24.244885
$ws.Range("N14").Value = 72.734655
$ws.Range("O14").Value = 0.6895205882382217
$ws.Range("P14").Value = 0.6895205882382218
$ws.Range("Q14").Value = 17.305683729495
$ws.Range("R14").Value = 155.751153565455
$ws.Range("S14").Value = 0.01820328093640073
$ws.Range("T14").Value = 0.01820328093640074

# Row 15
$ws.Range("A15").Value = "sCs"
$ws.Range("B15").Value = "Efna1"
$ws.Range("C15").Value = "Epha2"
$ws.Range("D15").Value = "FAPs"
$ws.Range("E15").Value = 3.0
$ws.Range("F15").Value = 1.0
$ws.Range("G15").Value = 0.713787
$ws.Range("H15").Value = 2.141361
$ws.Range("I15").Value = 0.02639990922230694
$ws.Range("J15").Value = 0.02639990922230694
$ws.Range("K15").Value = 3.0
$ws.Range("L15").Value = 1.0
$ws.Range("M15").Value = 0.7651789999999999
$ws.Range("N15").Value = 2.295537
$ws.Range("O15").Value = 0.02176156637523891
$ws.Range("P15").Value = 0.02176156637523891
$ws.Range("Q15").Value = 0.5461748228729999
$ws.Range("R15").Value = 4.915573405857
$ws.Range("S15").Value = 0.0005745033768415142
$ws.Range("T15").Value = 0.0005745033768415143

# Row 16
$ws.Range("A16").Value = "sCs"
$ws.Range("B16").Value = "Efna1"
$ws.Range("C16").Value = "Epha2"
$ws.Range("D16").Value = "M2"
$ws.Range("E16").Value = 3.0
$ws.Range("F16").Value = 1.0
$ws.Range("G16").Value = 0.713787
$ws.Range("H16").Value = 2.141361
$ws.Range("I16").Value = 0.02639990922230694
$ws.Range("J16").Value = 0.02639990922230694
$ws.Range("K16").Value = 3.0
$ws.Range("L16").Value = 1.0
$ws.Range("M16").Value = 0.6351283333333334
$ws.Range("N16").Value = 1.905385
$ws.Range("O16").Value = 0.01806294655581008
$ws.Range("P16").Value = 0.01806294655581008
$ws.Range("Q16").Value = 0.453346347665
$ws.Range("R16").Value = 4.080117128985
$ws.Range("S16").Value = 0.0004768601493607679
$ws.Range("T16").Value = 0.0004768601493607679

# Row 17
$ws.Range("A17").Value = "sCs"
$ws.Range("B17").Value = "Efna1"
$ws.Range("C17").Value = "Epha2"
$ws.Range("D17").Value = "sCs"
$ws.Range("E17").Value = 3.0
$ws.Range("F17").Value = 1.0
$ws.Range("G17").Value = 0.713787
$ws.Range("H17").Value = 2.141361
$ws.Range("I17").Value = 0.02639990922230694
$ws.Range("J17").Value = 0.02639990922230694
$ws.Range("K17").Value = 3.0
$ws.Range("L17").Value = 1.0
$ws.Range("M17").Value = 9.516752666666667
$ws.Range("N17").Value = 28.550258
$ws.Range("O17").Value = 0.2706548988307292
$ws.Range("P17").Value = 0.2706548988307293
$ws.Range("Q17").Value = 6.792934335682
$ws.Range("R17").Value = 61.13640902113799
$ws.Range("S17").Value = 0.007145264759703921
$ws.Range("T17").Value = 0.007145264759703922
